$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "YCHERN"
$ws.Range("C2").Value = "ASFLI"
$ws.Range("D2").Value = "REGISTERPROJECT"
$ws.Range("E2").Value = "PENDING"
$ws.Range("E2").Value = "APPROVED"
$ws.Range("F2").Value = 1
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Style = "Normal"
